$d = $word.ActiveDocument

# The ASSUNTO (subject) row of the header table currently reads:
#   "ATRIBUIÇÕES DA DIVISÃO DE NACIONALIZAÇÃO E QUALIFICAÇÃO"
# and must become:
#   "PROCESSOS DA DIVISÃO DE NACIONALIZAÇÃO E QUALIFICAÇÃO (NNAQ)"
# split across four runs: "PROCESSOS" / " DA DIVISÃO DE NACIONALIZAÇÃO E " /
# "QUALIFICAÇÃO" / " (NNAQ)".

# --- Step 1: swap "ATRIBUIÇÕES" for "PROCESSOS" ----------------------------
$rng = $d.Content
$rng.Find.Execute("ATRIBUIÇÕES", $true, $false, $false, $false, $false, $true, 1, $false, "PROCESSOS", 2) | Out-Null

# --- Step 2: append " (NNAQ)" right after "QUALIFICAÇÃO" -------------------
$rngAfterQual = $d.Content
$rngAfterQual.Find.Execute("QUALIFICAÇÃO") | Out-Null
$rngAfterQual.Collapse(0)                 # wdCollapseEnd
$rngAfterQual.InsertAfter(" (NNAQ)")

# --- Step 3: make sure each piece of text ends up in its own run -----------
# Toggling a character-formatting property on/off forces the run to be
# split at the range boundaries instead of being re-merged with its
# neighbour (the visible formatting is unchanged).
$rngProcessos = $d.Content
$rngProcessos.Find.Execute("PROCESSOS") | Out-Null
$rngProcessos.Font.Bold = 1
$rngProcessos.Font.Bold = 0

$rngQual = $d.Content
$rngQual.Find.Execute("QUALIFICAÇÃO") | Out-Null
$rngQual.Font.Bold = 1
$rngQual.Font.Bold = 0

$rngNNAQ = $d.Content
$rngNNAQ.Find.Execute(" (NNAQ)") | Out-Null
$rngNNAQ.Font.Bold = 1
$rngNNAQ.Font.Bold = 0
